# Generate Report for Handback
# Adds a new handback row (ba303ce5-72b3-4ebb-a135-ceb5b863046b.md) to the
# Overview, zh-cn and de-de sheets/tables of the handback-status workbook.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Sheet "Overview" (row 4)
# ---------------------------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")

$wsOverview.Range("A4").Value = "ba303ce5-72b3-4ebb-a135-ceb5b863046b.md"
$wsOverview.Hyperlinks.Add(
    $wsOverview.Range("B4"),
    "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/1f1a9b6b9d1f4a6c9e6b2a6f4d7c8e9a0b1c2d3e/e2e/ba303ce5-72b3-4ebb-a135-ceb5b863046b.md",
    "",
    "",
    "e2e\ba303ce5-72b3-4ebb-a135-ceb5b863046b.md"
)
$wsOverview.Range("C4").Value = ".md"
$wsOverview.Range("E4").Value = "Handed back: in sync with en-US"
$wsOverview.Range("F4").Value = "Handed back: in sync with en-US"
$wsOverview.Range("G4").Value = "2016-08-25 04:43:10"
$wsOverview.Range("G4").NumberFormat = "yyyy-mm-dd HH:mm:ss"

$loOverview = $wsOverview.ListObjects.Item(1)
$loOverview.Resize($wsOverview.Range("A1:G4"))

# ---------------------------------------------------------------------
# Sheet "zh-cn" (row 4)
# ---------------------------------------------------------------------
$wsZhCn = $wb.Worksheets.Item("zh-cn")

$wsZhCn.Hyperlinks.Add(
    $wsZhCn.Range("A4"),
    "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/1f1a9b6b9d1f4a6c9e6b2a6f4d7c8e9a0b1c2d3e/e2e/ba303ce5-72b3-4ebb-a135-ceb5b863046b.md",
    "",
    "",
    "ba303ce5-72b3-4ebb-a135-ceb5b863046b.md"
)
$wsZhCn.Range("B4").Value = ".md"
$wsZhCn.Range("C4").Value = "Handed back: in sync with en-US"
$wsZhCn.Range("D4").Value = "e2e"
$wsZhCn.Range("E4").Value = "ht"
$wsZhCn.Range("F4").Value = "'True"
$wsZhCn.Range("G4").Value = "ba303ce5-72b3-4ebb-a135-ceb5b863046b.560c0cb7ef31be98f43e29534b7421bf47740167.zh-cn.xlf"
$wsZhCn.Range("H4").Value = "2016-08-25 04:43:01"
$wsZhCn.Range("H4").NumberFormat = "yyyy-mm-dd HH:mm:ss"
$wsZhCn.Hyperlinks.Add(
    $wsZhCn.Range("I4"),
    "https://github.com/OpenLocalizationTestOrg/ol-test0-zhcn/blob/2a3b4c5d6e7f8091a2b3c4d5e6f708192a3b4c5d/e2e/ba303ce5-72b3-4ebb-a135-ceb5b863046b.md",
    "",
    "",
    "ba303ce5-72b3-4ebb-a135-ceb5b863046b.md"
)
$wsZhCn.Range("J4").Value = "ba303ce5-72b3-4ebb-a135-ceb5b863046b.560c0cb7ef31be98f43e29534b7421bf47740167.zh-cn.xlf"
$wsZhCn.Range("K4").Value = "2016-08-25 04:43:29"
$wsZhCn.Range("K4").NumberFormat = "yyyy-mm-dd HH:mm:ss"
$wsZhCn.Range("L4").Value = ""
$wsZhCn.Range("M4").Value = "'True"
$wsZhCn.Range("N4").Value = ""
$wsZhCn.Range("O4").Value = "'False"
$wsZhCn.Range("P4").Value = ""

$loZhCn = $wsZhCn.ListObjects.Item(1)
$loZhCn.Resize($wsZhCn.Range("A1:P4"))

# ---------------------------------------------------------------------
# Sheet "de-de" (row 4)
# ---------------------------------------------------------------------
$wsDeDe = $wb.Worksheets.Item("de-de")

$wsDeDe.Hyperlinks.Add(
    $wsDeDe.Range("A4"),
    "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/1f1a9b6b9d1f4a6c9e6b2a6f4d7c8e9a0b1c2d3e/e2e/ba303ce5-72b3-4ebb-a135-ceb5b863046b.md",
    "",
    "",
    "ba303ce5-72b3-4ebb-a135-ceb5b863046b.md"
)
$wsDeDe.Range("B4").Value = ".md"
$wsDeDe.Range("C4").Value = "Handed back: in sync with en-US"
$wsDeDe.Range("D4").Value = "e2e"
$wsDeDe.Range("E4").Value = "ht"
$wsDeDe.Range("F4").Value = "'True"
$wsDeDe.Range("G4").Value = "ba303ce5-72b3-4ebb-a135-ceb5b863046b.560c0cb7ef31be98f43e29534b7421bf47740167.de-de.xlf"
$wsDeDe.Range("H4").Value = "2016-08-25 04:43:10"
$wsDeDe.Range("H4").NumberFormat = "yyyy-mm-dd HH:mm:ss"
$wsDeDe.Hyperlinks.Add(
    $wsDeDe.Range("I4"),
    "https://github.com/OpenLocalizationTestOrg/ol-test0-dede/blob/3c4d5e6f708192a3b4c5d6e7f8091a2b3c4d5e6/e2e/ba303ce5-72b3-4ebb-a135-ceb5b863046b.md",
    "",
    "",
    "ba303ce5-72b3-4ebb-a135-ceb5b863046b.md"
)
$wsDeDe.Range("J4").Value = "ba303ce5-72b3-4ebb-a135-ceb5b863046b.560c0cb7ef31be98f43e29534b7421bf47740167.de-de.xlf"
$wsDeDe.Range("K4").Value = "2016-08-25 04:43:36"
$wsDeDe.Range("K4").NumberFormat = "yyyy-mm-dd HH:mm:ss"
$wsDeDe.Range("L4").Value = ""
$wsDeDe.Range("M4").Value = "'True"
$wsDeDe.Range("N4").Value = ""
$wsDeDe.Range("O4").Value = "'False"
$wsDeDe.Range("P4").Value = ""

$loDeDe = $wsDeDe.ListObjects.Item(1)
$loDeDe.Resize($wsDeDe.Range("A1:P4"))
